$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new entry (1 Euro -> TL question, euro.tlkur.com)
$ws.Range("A3").Value = 43158.08320918662
$ws.Range("A3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C3").Value = "1 Euro kaç Türk Lirası yapar?"
$ws.Range("I3").Value = "http://euro.tlkur.com"

# Row 4 - second visit of the same question, trailing-slash URL variant
$ws.Range("A4").Value = 43158.09216671144
$ws.Range("A4").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C4").Value = "1 Euro kaç Türk Lirası yapar?"
$ws.Range("I4").Value = "http://euro.tlkur.com/"
